# Bootstrap_Tables.docx update
# Refreshes the bootstrap estimates (B, SE, 95% CI, Sig) across the three
# tables to the new bootstrap run's results, and updates the CI-method
# footnote from "percentile" to "bca.simple" in all three table notes.

$d = $word.ActiveDocument

# ---- Table 1: Structural Path Coefficients From Bootstrap Analysis ----
$t = $d.Tables.Item(1)

# a1 (X -> M1)
$t.Cell(3,2).Range.Text = "0.127"
$t.Cell(3,3).Range.Text = "0.037"
$t.Cell(3,4).Range.Text = "[0.054, 0.200]"

# a1z (X x Z -> M1)
$t.Cell(4,2).Range.Text = "0.003"
$t.Cell(4,3).Range.Text = "0.015"
$t.Cell(4,4).Range.Text = "[-0.027, 0.032]"
$t.Cell(4,5).Range.Text = ""

# a2 (X -> M2)
$t.Cell(5,2).Range.Text = "-0.010"
$t.Cell(5,3).Range.Text = "0.036"
$t.Cell(5,4).Range.Text = "[-0.082, 0.061]"

# a2z (X x Z -> M2)
$t.Cell(6,2).Range.Text = "-0.014"
$t.Cell(6,3).Range.Text = "0.014"
$t.Cell(6,4).Range.Text = "[-0.042, 0.014]"
$t.Cell(6,5).Range.Text = ""

# c' (X -> Y)
$t.Cell(7,2).Range.Text = "0.041"
$t.Cell(7,3).Range.Text = "0.013"
$t.Cell(7,4).Range.Text = "[0.015, 0.068]"
$t.Cell(7,5).Range.Text = "*"

# c'z (X x Z -> Y)
$t.Cell(8,2).Range.Text = "-0.009"
$t.Cell(8,3).Range.Text = "0.005"
$t.Cell(8,4).Range.Text = "[-0.018, 0.000]"
$t.Cell(8,5).Range.Text = ""

# b1 (M1 -> Y)
$t.Cell(9,2).Range.Text = "-0.203"
$t.Cell(9,3).Range.Text = "0.008"
$t.Cell(9,4).Range.Text = "[-0.217, -0.188]"

# b2 (M2 -> Y)
$t.Cell(10,2).Range.Text = "0.160"
$t.Cell(10,3).Range.Text = "0.007"
$t.Cell(10,4).Range.Text = "[0.146, 0.174]"

# ---- Table 2: Direct, Indirect, and Total Effects at Mean Credit Dose ----
$t = $d.Tables.Item(2)

# Direct effect (c')
$t.Cell(3,2).Range.Text = "0.041"
$t.Cell(3,3).Range.Text = "0.013"
$t.Cell(3,4).Range.Text = "[0.015, 0.068]"
$t.Cell(3,5).Range.Text = "*"

# Indirect via EmoDiss
$t.Cell(4,2).Range.Text = "-0.026"
$t.Cell(4,3).Range.Text = "0.008"
$t.Cell(4,4).Range.Text = "[-0.041, -0.011]"

# Indirect via QualEngag
$t.Cell(5,2).Range.Text = "-0.002"
$t.Cell(5,3).Range.Text = "0.006"
$t.Cell(5,4).Range.Text = "[-0.013, 0.010]"

# Total effect
$t.Cell(6,2).Range.Text = "0.014"
$t.Cell(6,3).Range.Text = "0.017"
$t.Cell(6,4).Range.Text = "[-0.019, 0.047]"
$t.Cell(6,5).Range.Text = ""

# ---- Table 3: Conditional Indirect Effects and Index of Moderated Mediation ----
$t = $d.Tables.Item(3)

# Indirect via EmoDiss at -1 SD
$t.Cell(3,2).Range.Text = "-0.025"
$t.Cell(3,3).Range.Text = "0.012"
$t.Cell(3,4).Range.Text = "[-0.048, -0.001]"
$t.Cell(3,5).Range.Text = "*"

# Indirect via EmoDiss at Mean
$t.Cell(4,2).Range.Text = "-0.026"
$t.Cell(4,3).Range.Text = "0.008"
$t.Cell(4,4).Range.Text = "[-0.041, -0.011]"

# Indirect via EmoDiss at +1 SD
$t.Cell(5,2).Range.Text = "-0.027"
$t.Cell(5,3).Range.Text = "0.007"
$t.Cell(5,4).Range.Text = "[-0.040, -0.014]"

# IMM (EmoDiss)
$t.Cell(6,2).Range.Text = "-0.001"
$t.Cell(6,3).Range.Text = "0.003"
$t.Cell(6,4).Range.Text = "[-0.007, 0.005]"
$t.Cell(6,5).Range.Text = ""

# Indirect via QualEngag at -1 SD
$t.Cell(7,2).Range.Text = "0.003"
$t.Cell(7,3).Range.Text = "0.009"
$t.Cell(7,4).Range.Text = "[-0.015, 0.020]"
$t.Cell(7,5).Range.Text = ""

# Indirect via QualEngag at Mean
$t.Cell(8,2).Range.Text = "-0.002"
$t.Cell(8,3).Range.Text = "0.006"
$t.Cell(8,4).Range.Text = "[-0.013, 0.010]"

# Indirect via QualEngag at +1 SD
$t.Cell(9,2).Range.Text = "-0.006"
$t.Cell(9,3).Range.Text = "0.005"
$t.Cell(9,4).Range.Text = "[-0.016, 0.004]"
$t.Cell(9,5).Range.Text = ""

# IMM (QualEngag)
$t.Cell(10,2).Range.Text = "-0.002"
$t.Cell(10,3).Range.Text = "0.002"
$t.Cell(10,4).Range.Text = "[-0.007, 0.002]"
$t.Cell(10,5).Range.Text = ""

# ---- Footnotes: CI method name change (all three table notes) ----
$d.Content.Find.Execute("percentile CI excludes zero", $true, $false, $false, $false, $false, $true, 1, $false, "bca.simple CI excludes zero", 2)
